$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("E2").Value = "fecha,date"

# New "fecha" column values for existing rows, stored as text (Text number format)
$ws.Range("E3:E7").NumberFormat = "@"
$ws.Range("E3").Value = "1991-12-12"
$ws.Range("E4").Value = "1992-05-05"
$ws.Range("E5").Value = "1993-04-03"

# New rows 6 and 7
$ws.Range("A6").Value = "Jose"
$ws.Range("B6").Value = "Perezz"
$ws.Range("C6").Value = 414
$ws.Range("D6").Value = "m"
$ws.Range("E6").Value = "1993-02-03"

$ws.Range("A7").Value = "Jose"
$ws.Range("B7").Value = "Perezz"
$ws.Range("C7").Value = 414
$ws.Range("D7").Value = "m"
$ws.Range("E7").Value = "1993-02-03"

$ws.Range("A7:E7").Select()
